# Atualização automática de preços de eletricidade
# Updates the spot price data row (row 2) in the Spot_PT sheet with the
# latest day's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45962
$ws.Range("B2").Value = 45.93
$ws.Range("C2").Value = 24.42
$ws.Range("D2").Value = 22.49
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 30.09
$ws.Range("G2").Value = 33.86
$ws.Range("H2").Value = 47.3
$ws.Range("I2").Value = 52
$ws.Range("J2").Value = 50.4
$ws.Range("K2").Value = 23.75
$ws.Range("L2").Value = 18.3
$ws.Range("M2").Value = 14.9
$ws.Range("N2").Value = 10.78
$ws.Range("O2").Value = 15.65
$ws.Range("P2").Value = 15.21
$ws.Range("Q2").Value = 21.01
$ws.Range("R2").Value = 36.22
$ws.Range("S2").Value = 56.66
$ws.Range("T2").Value = 82.40000000000001
$ws.Range("U2").Value = 103.38
$ws.Range("V2").Value = 104.94
$ws.Range("W2").Value = 100.43
$ws.Range("X2").Value = 80.98999999999999
$ws.Range("Y2").Value = 72.87
$ws.Range("Z2").Value = 45.58
$ws.Range("AB2").Value = 89.81
$ws.Range("AD2").Value = 102.68
$ws.Range("AF2").Value = 92.89
